$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D2').Value = '62.119.80'
$ws.Range('D3').Value = '3.432.00'
$ws.Range('D5').Value = "'406.83"
$ws.Range('D6').Value = "'132.85"
$ws.Range('D7').Value = "'0.594"
$ws.Range('D9').Value = "'0.687"
$ws.Range('D10').Value = "'0.127"
$ws.Range('D11').Value = "'42.00"
$ws.Range('D13').Value = "'8.47"
$ws.Range('D14').Value = "'19.85"
$ws.Range('D15').Value = '3.433.17'
$ws.Range('D16').Value = '62.213.80'
$ws.Range('D18').Value = "'1.02"
$ws.Range('D19').Value = "'0.0000142"
$ws.Range('D20').Value = "'3.18"
$ws.Range('D21').Value = "'83.67"
$ws.Range('D22').Value = "'310.85"
$ws.Range('D23').Value = "'12.86"
$ws.Range('D24').Value = "'3.15"
$ws.Range('D25').Value = "'4.77"
$ws.Range('D26').Value = "'29.73"
$ws.Range('D27').Value = "'8.12"
$ws.Range('D28').Value = "'7.68"
$ws.Range('D29').Value = "'2.77"
$ws.Range('D30').Value = "'0.173"
$ws.Range('D32').Value = "'42.82"
$ws.Range('D34').Value = "'11.37"
$ws.Range('D35').Value = "'0.0486"
$ws.Range('D36').Value = "'51.41"
$ws.Range('D37').Value = "'1.00"
$ws.Range('D38').Value = "'3.39"
$ws.Range('D39').Value = "'0.326"
$ws.Range('D40').Value = "'2.92"
$ws.Range('D41').Value = "'138.44"
$ws.Range('D42').Value = "'1.98"
$ws.Range('D44').Value = "'3.98"
$ws.Range('D45').Value = "'16.74"
$ws.Range('D47').Value = "'21.12"
$ws.Range('D48').Value = '2.121.39'
$ws.Range('D50').Value = "'1.75"
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('E11').Value = '  -3.70%  '
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('E13').Value = '  -3.97%  '
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('E17').Value = '  +4.91%  '
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('E20').Value = '  -5.52%  '
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('E25').Value = '  +8.72%  '
$ws.Range('E26').Value = '  -3.06%  '
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('E29').Value = '  +4.61%  '
$ws.Range('E31').Value = '  -3.74%  '
$ws.Range('E32').Value = '  -3.22%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('E38').Value = '  -6.28%  '
$ws.Range('E39').Value = '  +12.67%  '
$ws.Range('E40').Value = '  -4.15%  '
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('E45').Value = '  -4.16%  '
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('E47').Value = '  -5.68%  '
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('E49').Value = '  -3.69%  '
$ws.Range('E50').Value = '  +20.93%  '
$ws.Range('E51').Value = '  +1.58%  '
